# Updates cryptos list figures (price / 1h volume change) to the latest
# scrape, matching the GitHub Actions commit "Updated cryptos list".
# A handful of rows were also re-ranked (coin name/link/price/%% swapped
# with a neighboring row).
#
# Cell values are plain text (e.g. "41.463.47", "  -2.20%  ") in the
# source workbook, so we force the NumberFormat to Text before writing
# each value -- otherwise Excel auto-coerces anything that looks like a
# number (e.g. "309.70") into a float and mangles it (trailing zeros,
# floating point noise, scientific notation, ...). Resetting the style
# to "Normal" afterwards keeps the cell formatting identical to the
# original (no explicit style / default style), it only affects the
# stored value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.463.47'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.20%  '
$ws.Range('E2').Style = 'Normal'
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.439.25'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.72%  '
$ws.Range('E3').Style = 'Normal'
# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +1.31%  '
$ws.Range('E4').Style = 'Normal'
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('E5').Style = 'Normal'
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '89.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -7.97%  '
$ws.Range('E6').Style = 'Normal'
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.532'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -5.05%  '
$ws.Range('E7').Style = 'Normal'
# Row 8
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.18%  '
$ws.Range('E8').Style = 'Normal'
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.480'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -7.29%  '
$ws.Range('E9').Style = 'Normal'
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.01'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -8.01%  '
$ws.Range('E10').Style = 'Normal'
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0763'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.20%  '
$ws.Range('E11').Style = 'Normal'
# Row 12
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('E12').Style = 'Normal'
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.824.31'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.45%  '
$ws.Range('E13').Style = 'Normal'
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -6.27%  '
$ws.Range('E14').Style = 'Normal'
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.512.23'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('E15').Style = 'Normal'
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.60%  '
$ws.Range('E16').Style = 'Normal'
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.760'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.60%  '
$ws.Range('E17').Style = 'Normal'
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.142.51'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.07%  '
$ws.Range('E18').Style = 'Normal'
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.14'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -5.66%  '
$ws.Range('E19').Style = 'Normal'
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0897'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -3.79%  '
$ws.Range('E20').Style = 'Normal'
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('E21').Style = 'Normal'
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.70'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -10.64%  '
$ws.Range('E22').Style = 'Normal'
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '232.09'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.47%  '
$ws.Range('E23').Style = 'Normal'
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.66'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.32%  '
$ws.Range('E24').Style = 'Normal'
# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.13%  '
$ws.Range('E25').Style = 'Normal'
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.83'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -6.63%  '
$ws.Range('E26').Style = 'Normal'
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.41'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -6.97%  '
$ws.Range('E27').Style = 'Normal'
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.40%  '
$ws.Range('E28').Style = 'Normal'
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.60%  '
$ws.Range('E29').Style = 'Normal'
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.75'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -6.90%  '
$ws.Range('E30').Style = 'Normal'
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '151.72'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.45%  '
$ws.Range('E31').Style = 'Normal'
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.27'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -9.19%  '
$ws.Range('E32').Style = 'Normal'
# Row 33
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.54'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -3.23%  '
$ws.Range('E33').Style = 'Normal'
# Row 34
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'ApeXProtocol'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.54'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.16%  '
$ws.Range('E34').Style = 'Normal'
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0736'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -5.11%  '
$ws.Range('E35').Style = 'Normal'
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '16.94'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.87%  '
$ws.Range('E36').Style = 'Normal'
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.91'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -5.62%  '
$ws.Range('E37').Style = 'Normal'
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.78'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -8.49%  '
$ws.Range('E38').Style = 'Normal'
# Row 39
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -4.18%  '
$ws.Range('E39').Style = 'Normal'
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0978'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -9.53%  '
$ws.Range('E40').Style = 'Normal'
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.90'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -6.36%  '
$ws.Range('E41').Style = 'Normal'
# Row 42
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '20.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.99%  '
$ws.Range('E42').Style = 'Normal'
# Row 43
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.02'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.51%  '
$ws.Range('E43').Style = 'Normal'
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.928.73'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.61%  '
$ws.Range('E44').Style = 'Normal'
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0274'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -6.17%  '
$ws.Range('E45').Style = 'Normal'
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.86'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -9.68%  '
$ws.Range('E46').Style = 'Normal'
# Row 47
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.52'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.12%  '
$ws.Range('E47').Style = 'Normal'
# Row 48
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.701.58'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.07%  '
$ws.Range('E48').Style = 'Normal'
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '93.67'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.61%  '
$ws.Range('E49').Style = 'Normal'
# Row 50
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Algorand'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.171'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -7.79%  '
$ws.Range('E50').Style = 'Normal'
# Row 51
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.31'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -9.54%  '
$ws.Range('E51').Style = 'Normal'
